# Update the worker "database" rows in the account statement (Estado de Cuenta).
# Rows 16 and 17 list workers (doc type, doc number, name, period) together with
# their "Valor Mora" (overdue amount). The edit swaps the two workers' identity
# data between the rows while each worker keeps their own overdue amount.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values for row 16 and row 17 (doc number, name, overdue value)
$docNum16 = $ws.Range("C16").Value()
$name16   = $ws.Range("D16").Value()
$valor16  = $ws.Range("G16").Value()

$docNum17 = $ws.Range("C17").Value()
$name17   = $ws.Range("D17").Value()
$valor17  = $ws.Range("G17").Value()

# Swap them: worker previously on row 17 moves to row 16 (keeping their own
# overdue value attached), and vice versa.
$ws.Range("C16").Value = $docNum17
$ws.Range("D16").Value = $name17
$ws.Range("G16").Value = $valor17

$ws.Range("C17").Value = $docNum16
$ws.Range("D17").Value = $name16
$ws.Range("G17").Value = $valor16
